$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 998
$ws.Range("J19").Value = 1140.6666
$ws.Range("L19").Value = 1140.6666
$ws.Range("N19").Value = -1490.6666
$ws.Range("H40").Value = 8581.666999999999
$ws.Range("J40").Value = 9107.666999999999
$ws.Range("L40").Value = 9107.666999999999
$ws.Range("N40").Value = -9457.666999999999
$ws.Range("H70").Value = 5744.5
$ws.Range("I70").Value = 3940
$ws.Range("J70").Value = 6195.625
$ws.Range("K70").Value = 11820
$ws.Range("L70").Value = 18586.875
$ws.Range("M70").Value = -11550
$ws.Range("N70").Value = -19126.875
$ws.Range("H73").Value = 5744.5
$ws.Range("I73").Value = 3940
$ws.Range("J73").Value = 6195.625
$ws.Range("K73").Value = 11820
$ws.Range("L73").Value = 18586.875
$ws.Range("M73").Value = -10884
$ws.Range("N73").Value = -20458.875
$ws.Range("H80").Value = 1307.1
$ws.Range("I80").Value = 407.75
$ws.Range("J80").Value = 1906.6666
$ws.Range("K80").Value = 1223.25
$ws.Range("L80").Value = 5719.9998
$ws.Range("M80").Value = -225.25
$ws.Range("N80").Value = -7715.9998
$ws.Range("H83").Value = 1307.1
$ws.Range("I83").Value = 407.75
$ws.Range("J83").Value = 1906.6666
$ws.Range("K83").Value = 3669.75
$ws.Range("L83").Value = 17159.9994
$ws.Range("M83").Value = 1322.25
$ws.Range("N83").Value = -27143.9994
$ws.Range("H106").Value = 88004090
$ws.Range("I106").Value = 110003450
$ws.Range("J106").Value = 6659
$ws.Range("K106").Value = 110003450
$ws.Range("L106").Value = 6659
$ws.Range("M106").Value = -110002819
$ws.Range("N106").Value = -7921
$ws.Range("H113").Value = 2176704
$ws.Range("I113").Value = 5558068
$ws.Range("J113").Value = 2969.9285
$ws.Range("K113").Value = 5558068
$ws.Range("L113").Value = 2969.9285
$ws.Range("M113").Value = -5554814
$ws.Range("N113").Value = -9477.9285
$ws.Range("H132").Value = 3419.681
$ws.Range("I132").Value = 2249.4187
$ws.Range("K132").Value = 6748.256100000001
$ws.Range("M132").Value = -4218.256100000001
$ws.Range("H135").Value = 448.0435
$ws.Range("I135").Value = 432.75
$ws.Range("K135").Value = 3894.75
$ws.Range("M135").Value = -1359.75
$ws.Range("H137").Value = 2109.5789
$ws.Range("I137").Value = 1853.2727
$ws.Range("J137").Value = 2462
$ws.Range("K137").Value = 5559.8181
$ws.Range("L137").Value = 7386
$ws.Range("M137").Value = -3009.8181
$ws.Range("N137").Value = -12486
$ws.Range("H138").Value = 2331.0652
$ws.Range("I138").Value = 2365.75
$ws.Range("J138").Value = 2312.5667
$ws.Range("K138").Value = 7097.25
$ws.Range("L138").Value = 6937.7001
$ws.Range("M138").Value = -1957.25
$ws.Range("N138").Value = -17217.7001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7918
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 7918
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 7918
$ws.Range("M74").Value = $null
$ws.Range("N74").Value = -9666
$ws.Range("H77").Value = 7918
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 7918
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 39590
$ws.Range("M77").Value = $null
$ws.Range("N77").Value = -48326
$ws.Range("H132").Value = 2771.7
$ws.Range("I132").Value = 1714.6875
$ws.Range("J132").Value = 6999.75
$ws.Range("K132").Value = 5144.0625
$ws.Range("L132").Value = 20999.25
$ws.Range("M132").Value = -2614.0625
$ws.Range("N132").Value = -26059.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3055.6667
$ws.Range("I134").Value = 2552.16
$ws.Range("J134").Value = 9349.5
$ws.Range("K134").Value = 7656.48
$ws.Range("L134").Value = 28048.5
$ws.Range("M134").Value = -5121.48
$ws.Range("N134").Value = -33118.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3434.111
$ws.Range("I58").Value = 3703
$ws.Range("J58").Value = 3263
$ws.Range("K58").Value = 3703
$ws.Range("L58").Value = 3263
$ws.Range("M58").Value = -3500
$ws.Range("N58").Value = -3669
$ws.Range("H62").Value = 7937.1665
$ws.Range("J62").Value = 8370.308000000001
$ws.Range("L62").Value = 8370.308000000001
$ws.Range("N62").Value = -9618.308000000001
$ws.Range("H65").Value = 7937.1665
$ws.Range("J65").Value = 8370.308000000001
$ws.Range("L65").Value = 41851.54000000001
$ws.Range("N65").Value = -48091.54000000001
$ws.Range("H136").Value = 3434.111
$ws.Range("I136").Value = 3703
$ws.Range("J136").Value = 3263
$ws.Range("K136").Value = 11109
$ws.Range("L136").Value = 9789
$ws.Range("M136").Value = -8559
$ws.Range("N136").Value = -14889

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 443.75
$ws.Range("I47").Value = 443.66666
$ws.Range("K47").Value = 1330.99998
$ws.Range("M47").Value = -899.9999800000001
$ws.Range("H98").Value = 576.625
$ws.Range("I98").Value = 332.66666
$ws.Range("J98").Value = 723
$ws.Range("K98").Value = 997.9999799999999
$ws.Range("L98").Value = 2169
$ws.Range("M98").Value = 500.0000200000001
$ws.Range("N98").Value = -5165
$ws.Range("H113").Value = 701.4545000000001
$ws.Range("J113").Value = 912.53845
$ws.Range("L113").Value = 2737.61535
$ws.Range("N113").Value = -7077.61535
$ws.Range("H136").Value = 2131.2727
$ws.Range("I136").Value = 2131.2727
$ws.Range("K136").Value = 6393.8181
$ws.Range("M136").Value = -1293.8181
$ws.Range("H137").Value = 3492.5334
$ws.Range("I137").Value = 2148.2222
$ws.Range("J137").Value = 5509
$ws.Range("K137").Value = 6444.6666
$ws.Range("L137").Value = 16527
$ws.Range("M137").Value = -1344.6666
$ws.Range("N137").Value = -26727
$ws.Range("H138").Value = 2247.8125
$ws.Range("I138").Value = 1955.3572
$ws.Range("J138").Value = 4295
$ws.Range("K138").Value = 5866.071599999999
$ws.Range("L138").Value = 12885
$ws.Range("M138").Value = -726.0715999999993
$ws.Range("N138").Value = -23165
$ws.Range("H139").Value = 2832
$ws.Range("I139").Value = 2850.8572
$ws.Range("J139").Value = 2700
$ws.Range("K139").Value = 8552.571599999999
$ws.Range("L139").Value = 8100
$ws.Range("M139").Value = -3412.571599999999
$ws.Range("N139").Value = -18380
$ws.Range("H140").Value = 1830.0526
$ws.Range("I140").Value = 1302.6
$ws.Range("K140").Value = 3907.8
$ws.Range("M140").Value = 1272.2
$ws.Range("H141").Value = 3227.1538
$ws.Range("I141").Value = 3079.4167
$ws.Range("K141").Value = 9238.250100000001
$ws.Range("M141").Value = -4058.250100000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2494.6428
$ws.Range("I80").Value = 2515.9092
$ws.Range("J80").Value = 2416.6667
$ws.Range("K80").Value = 2515.9092
$ws.Range("L80").Value = 2416.6667
$ws.Range("M80").Value = -1517.9092
$ws.Range("N80").Value = -4412.6667
$ws.Range("H83").Value = 2494.6428
$ws.Range("I83").Value = 2515.9092
$ws.Range("J83").Value = 2416.6667
$ws.Range("K83").Value = 12579.546
$ws.Range("L83").Value = 12083.3335
$ws.Range("M83").Value = -7587.546
$ws.Range("N83").Value = -22067.3335
$ws.Range("H132").Value = 4403.7256
$ws.Range("I132").Value = 2994.353
$ws.Range("J132").Value = 7222.4707
$ws.Range("K132").Value = 8983.059000000001
$ws.Range("L132").Value = 21667.4121
$ws.Range("M132").Value = -6453.059000000001
$ws.Range("N132").Value = -26727.4121

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4409.6
$ws.Range("I40").Value = 2243.6875
$ws.Range("J40").Value = 8260.111000000001
$ws.Range("K40").Value = 2243.6875
$ws.Range("L40").Value = 8260.111000000001
$ws.Range("M40").Value = -2107.6875
$ws.Range("N40").Value = -8532.111000000001
$ws.Range("H82").Value = 2274.3
$ws.Range("I82").Value = 1416.4
$ws.Range("J82").Value = 3132.2
$ws.Range("K82").Value = 1416.4
$ws.Range("L82").Value = 3132.2
$ws.Range("M82").Value = -1055.4
$ws.Range("N82").Value = -3854.2
$ws.Range("H85").Value = 2274.3
$ws.Range("I85").Value = 1416.4
$ws.Range("J85").Value = 3132.2
$ws.Range("K85").Value = 1416.4
$ws.Range("L85").Value = 3132.2
$ws.Range("M85").Value = -168.4000000000001
$ws.Range("N85").Value = -5628.2
$ws.Range("H136").Value = 7476.4165
$ws.Range("I136").Value = 6509.2583
$ws.Range("K136").Value = 19527.7749
$ws.Range("M136").Value = -16977.7749

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3816.8262
$ws.Range("I136").Value = 2559.1064
$ws.Range("J136").Value = 6503.773
$ws.Range("K136").Value = 7677.3192
$ws.Range("L136").Value = 19511.319
$ws.Range("M136").Value = -5127.3192
$ws.Range("N136").Value = -24611.319
